# Shift the tick-price time series forward: each row now holds the values
# that used to live 17 rows further down, and 4 brand-new trading days are
# appended at the end; the 13 oldest trailing rows are dropped so the sheet
# shrinks from A1:B269 down to A1:B256.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{r=2; d=45051; p="31.760"},
    @{r=3; d=45054; p="31.970"},
    @{r=4; d=45055; p="31.990"},
    @{r=5; d=45056; p="32.040"},
    @{r=6; d=45057; p="31.680"},
    @{r=7; d=45058; p="31.230"},
    @{r=8; d=45061; p="30.760"},
    @{r=9; d=45062; p="30.940"},
    @{r=10; d=45063; p="30.990"},
    @{r=11; d=45064; p="31.220"},
    @{r=12; d=45065; p="31.280"},
    @{r=13; d=45068; p="31.490"},
    @{r=14; d=45069; p="31.090"},
    @{r=15; d=45070; p="30.820"},
    @{r=16; d=45071; p="30.930"},
    @{r=17; d=45072; p="31.110"},
    @{r=18; d=45075; p="31.180"},
    @{r=19; d=45076; p="31.270"},
    @{r=20; d=45077; p="31.120"},
    @{r=21; d=45078; p="31.560"},
    @{r=22; d=45079; p="31.320"},
    @{r=23; d=45082; p="31.400"},
    @{r=24; d=45083; p="31.250"},
    @{r=25; d=45084; p="31.160"},
    @{r=26; d=45085; p="31.310"},
    @{r=27; d=45086; p="31.480"},
    @{r=28; d=45089; p="31.870"},
    @{r=29; d=45090; p="31.380"},
    @{r=30; d=45091; p="31.230"},
    @{r=31; d=45092; p="29.829"},
    @{r=32; d=45093; p="30.124"},
    @{r=33; d=45096; p="29.648"},
    @{r=34; d=45097; p="29.629"},
    @{r=35; d=45098; p="29.771"},
    @{r=36; d=45099; p="29.610"},
    @{r=37; d=45100; p="29.495"},
    @{r=38; d=45103; p="29.552"},
    @{r=39; d=45104; p="29.800"},
    @{r=40; d=45105; p="29.943"},
    @{r=41; d=45106; p="30.200"},
    @{r=42; d=45107; p="30.638"},
    @{r=43; d=45110; p="30.257"},
    @{r=44; d=45111; p="30.120"},
    @{r=45; d=45112; p="29.810"},
    @{r=46; d=45113; p="29.230"},
    @{r=47; d=45114; p="29.420"},
    @{r=48; d=45117; p="29.610"},
    @{r=49; d=45118; p="30.050"},
    @{r=50; d=45119; p="30.160"},
    @{r=51; d=45120; p="29.980"},
    @{r=52; d=45121; p="29.760"},
    @{r=53; d=45124; p="29.640"},
    @{r=54; d=45125; p="29.740"},
    @{r=55; d=45126; p="30.300"},
    @{r=56; d=45127; p="30.800"},
    @{r=57; d=45128; p="30.660"},
    @{r=58; d=45131; p="30.450"},
    @{r=59; d=45132; p="30.580"},
    @{r=60; d=45133; p="31.550"},
    @{r=61; d=45134; p="31.840"},
    @{r=62; d=45135; p="31.940"},
    @{r=63; d=45138; p="31.790"},
    @{r=64; d=45139; p="31.640"},
    @{r=65; d=45140; p="31.430"},
    @{r=66; d=45141; p="31.150"},
    @{r=67; d=45142; p="31.360"},
    @{r=68; d=45145; p="31.440"},
    @{r=69; d=45146; p="31.390"},
    @{r=70; d=45147; p="31.640"},
    @{r=71; d=45148; p="31.770"},
    @{r=72; d=45149; p="31.800"},
    @{r=73; d=45152; p="31.680"},
    @{r=74; d=45153; p="31.510"},
    @{r=75; d=45154; p="31.530"},
    @{r=76; d=45155; p="31.280"},
    @{r=77; d=45156; p="31.300"},
    @{r=78; d=45159; p="31.310"},
    @{r=79; d=45160; p="31.740"},
    @{r=80; d=45161; p="31.750"},
    @{r=81; d=45162; p="31.920"},
    @{r=82; d=45163; p="31.670"},
    @{r=83; d=45166; p="32.020"},
    @{r=84; d=45167; p="32.180"},
    @{r=85; d=45168; p="32.340"},
    @{r=86; d=45169; p="32.390"},
    @{r=87; d=45170; p="32.330"},
    @{r=88; d=45173; p="32.370"},
    @{r=89; d=45174; p="32.260"},
    @{r=90; d=45175; p="32.240"},
    @{r=91; d=45176; p="32.390"},
    @{r=92; d=45177; p="32.480"},
    @{r=93; d=45180; p="32.940"},
    @{r=94; d=45181; p="32.860"},
    @{r=95; d=45182; p="32.660"},
    @{r=96; d=45183; p="32.940"},
    @{r=97; d=45184; p="33.410"},
    @{r=98; d=45187; p="33.620"},
    @{r=99; d=45188; p="34.250"},
    @{r=100; d=45189; p="34.360"},
    @{r=101; d=45190; p="34.160"},
    @{r=102; d=45191; p="34.060"},
    @{r=103; d=45194; p="33.830"},
    @{r=104; d=45195; p="33.880"},
    @{r=105; d=45196; p="33.800"},
    @{r=106; d=45197; p="33.940"},
    @{r=107; d=45198; p="34.070"},
    @{r=108; d=45201; p="33.720"},
    @{r=109; d=45202; p="33.450"},
    @{r=110; d=45203; p="33.360"},
    @{r=111; d=45204; p="32.030"},
    @{r=112; d=45205; p="32.920"},
    @{r=113; d=45208; p="32.890"},
    @{r=114; d=45209; p="33.440"},
    @{r=115; d=45210; p="33.530"},
    @{r=116; d=45211; p="33.620"},
    @{r=117; d=45212; p="33.200"},
    @{r=118; d=45215; p="33.130"},
    @{r=119; d=45216; p="33.110"},
    @{r=120; d=45217; p="32.750"},
    @{r=121; d=45218; p="32.510"},
    @{r=122; d=45219; p="32.250"},
    @{r=123; d=45222; p="32.410"},
    @{r=124; d=45223; p="32.620"},
    @{r=125; d=45224; p="32.510"},
    @{r=126; d=45225; p="32.750"},
    @{r=127; d=45226; p="33.420"},
    @{r=128; d=45229; p="33.740"},
    @{r=129; d=45230; p="34.170"},
    @{r=130; d=45231; p="33.990"},
    @{r=131; d=45232; p="34.190"},
    @{r=132; d=45233; p="33.890"},
    @{r=133; d=45236; p="33.500"},
    @{r=134; d=45237; p="33.180"},
    @{r=135; d=45238; p="33.260"},
    @{r=136; d=45239; p="33.320"},
    @{r=137; d=45240; p="33.260"},
    @{r=138; d=45243; p="33.780"},
    @{r=139; d=45244; p="34.430"},
    @{r=140; d=45245; p="35.190"},
    @{r=141; d=45246; p="35.560"},
    @{r=142; d=45247; p="35.890"},
    @{r=143; d=45250; p="35.960"},
    @{r=144; d=45251; p="36.180"},
    @{r=145; d=45252; p="36.060"},
    @{r=146; d=45253; p="36.000"},
    @{r=147; d=45254; p="36.230"},
    @{r=148; d=45257; p="36.410"},
    @{r=149; d=45258; p="36.370"},
    @{r=150; d=45259; p="36.760"},
    @{r=151; d=45260; p="36.650"},
    @{r=152; d=45261; p="37.050"},
    @{r=153; d=45264; p="36.930"},
    @{r=154; d=45265; p="37.170"},
    @{r=155; d=45266; p="37.070"},
    @{r=156; d=45267; p="37.060"},
    @{r=157; d=45268; p="37.450"},
    @{r=158; d=45271; p="37.460"},
    @{r=159; d=45272; p="37.980"},
    @{r=160; d=45273; p="39.070"},
    @{r=161; d=45274; p="39.580"},
    @{r=162; d=45275; p="39.780"},
    @{r=163; d=45278; p="39.300"},
    @{r=164; d=45279; p="39.350"},
    @{r=165; d=45280; p="39.480"},
    @{r=166; d=45281; p="39.770"},
    @{r=167; d=45282; p="39.870"},
    @{r=168; d=45287; p="40.160"},
    @{r=169; d=45288; p="40.110"},
    @{r=170; d=45289; p="40.160"},
    @{r=171; d=45293; p="40.350"},
    @{r=172; d=45294; p="39.410"},
    @{r=173; d=45295; p="39.940"},
    @{r=174; d=45296; p="39.700"},
    @{r=175; d=45299; p="39.980"},
    @{r=176; d=45300; p="39.540"},
    @{r=177; d=45301; p="38.977"},
    @{r=178; d=45302; p="39.125"},
    @{r=179; d=45303; p="39.728"},
    @{r=180; d=45306; p="39.797"},
    @{r=181; d=45307; p="39.313"},
    @{r=182; d=45308; p="39.224"},
    @{r=183; d=45309; p="39.214"},
    @{r=184; d=45310; p="39.440"},
    @{r=185; d=45313; p="39.750"},
    @{r=186; d=45314; p="39.330"},
    @{r=187; d=45315; p="39.660"},
    @{r=188; d=45316; p="39.500"},
    @{r=189; d=45317; p="39.620"},
    @{r=190; d=45320; p="35.660"},
    @{r=191; d=45321; p="36.240"},
    @{r=192; d=45322; p="36.590"},
    @{r=193; d=45323; p="36.610"},
    @{r=194; d=45324; p="36.860"},
    @{r=195; d=45327; p="36.460"},
    @{r=196; d=45328; p="36.500"},
    @{r=197; d=45329; p="36.050"},
    @{r=198; d=45330; p="35.940"},
    @{r=199; d=45331; p="35.730"},
    @{r=200; d=45334; p="36.210"},
    @{r=201; d=45335; p="36.390"},
    @{r=202; d=45336; p="36.510"},
    @{r=203; d=45337; p="36.930"},
    @{r=204; d=45338; p="36.390"},
    @{r=205; d=45341; p="36.450"},
    @{r=206; d=45342; p="36.630"},
    @{r=207; d=45343; p="36.730"},
    @{r=208; d=45344; p="37.070"},
    @{r=209; d=45345; p="37.160"},
    @{r=210; d=45348; p="36.940"},
    @{r=211; d=45349; p="37.530"},
    @{r=212; d=45350; p="37.770"},
    @{r=213; d=45351; p="37.940"},
    @{r=214; d=45352; p="37.130"},
    @{r=215; d=45355; p="37.740"},
    @{r=216; d=45356; p="38.000"},
    @{r=217; d=45357; p="38.500"},
    @{r=218; d=45358; p="39.150"},
    @{r=219; d=45359; p="39.070"},
    @{r=220; d=45362; p="38.770"},
    @{r=221; d=45363; p="39.040"},
    @{r=222; d=45364; p="39.930"},
    @{r=223; d=45365; p="40.120"},
    @{r=224; d=45366; p="39.860"},
    @{r=225; d=45369; p="40.140"},
    @{r=226; d=45370; p="39.700"},
    @{r=227; d=45371; p="40.350"},
    @{r=228; d=45372; p="40.590"},
    @{r=229; d=45373; p="40.460"},
    @{r=230; d=45376; p="40.070"},
    @{r=231; d=45377; p="40.000"},
    @{r=232; d=45378; p="40.040"},
    @{r=233; d=45379; p="38.780"},
    @{r=234; d=45384; p="38.480"},
    @{r=235; d=45385; p="38.720"},
    @{r=236; d=45386; p="38.240"},
    @{r=237; d=45387; p="37.900"},
    @{r=238; d=45390; p="37.880"},
    @{r=239; d=45391; p="37.780"},
    @{r=240; d=45392; p="37.540"},
    @{r=241; d=45393; p="37.600"},
    @{r=242; d=45394; p="37.920"},
    @{r=243; d=45397; p="37.640"},
    @{r=244; d=45398; p="38.140"},
    @{r=245; d=45399; p="38.660"},
    @{r=246; d=45400; p="38.280"},
    @{r=247; d=45401; p="38.220"},
    @{r=248; d=45404; p="38.460"},
    @{r=249; d=45405; p="38.480"},
    @{r=250; d=45406; p="38.360"},
    @{r=251; d=45407; p="37.700"},
    @{r=252; d=45408; p="37.820"},
    @{r=253; d=45411; p="37.980"},
    @{r=254; d=45412; p="37.580"},
    @{r=255; d=45414; p="37.800"},
    @{r=256; d=45415; p="37.860"}
)

foreach ($item in $rows) {
    if ($item.r -ge 253) {
        $ws.Cells.Item($item.r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
    $ws.Cells.Item($item.r, 1).Value = $item.d
    $ws.Cells.Item($item.r, 2).NumberFormat = "@"
    $ws.Cells.Item($item.r, 2).Value = $item.p
    $ws.Cells.Item($item.r, 2).ClearFormats()
}

# Drop the 13 rows that fell off the end of the window.
$ws.Range("A257:B269").EntireRow.Delete()
